# Predicted_LaLiga2025_26_table_matchday_9.xlsx
# Add WIN / TOP4 / TOP5 / TOP6 / RELEGATION placeholder columns (C:G),
# move the ExpPoints column to H, and refresh the team order/points
# for the new matchday-9 prediction (Monte-Carlo simulation prep).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "WIN"
$ws.Cells.Item(1, 4).Value = "TOP4"
$ws.Cells.Item(1, 5).Value = "TOP5"
$ws.Cells.Item(1, 6).Value = "TOP6"
$ws.Cells.Item(1, 7).Value = "RELEGATION"
$ws.Cells.Item(1, 8).Value = "ExpPoints"

# Copy the header formatting (style s="1") from an existing header cell
# onto the newly added header cells D1:H1.
$ws.Range("A1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows ----------------------------------------------------------
# row, rank, team, expPoints(new, col H)
$data = @(
    @(2,  1, "Real Madrid",          89.56181065806307),
    @(3,  2, "Barcelona",            88.59847473929761),
    @(4,  3, "Atlético de Madrid",   70.78324552510308),
    @(5,  4, "Villarreal",           65.75200274297808),
    @(6,  5, "Real Betis",           60.24138812159953),
    @(7,  6, "Athletic Club",        58.5950811416851),
    @(8,  7, "Rayo Vallecano",       54.79762064968454),
    @(9,  8, "Sevilla",              48.93833423376397),
    @(10, 9, "Getafe",               48.38458984467765),
    @(11, 10, "Celta de Vigo",       48.11360744718921),
    @(12, 11, "Osasuna",             47.71113756474548),
    @(13, 12, "Valencia",            46.68502906762888),
    @(14, 13, "Real Sociedad",       46.20594350989916),
    @(15, 14, "Espanyol",            45.81938575044128),
    @(16, 15, "Alavés",              39.66044537339634),
    @(17, 16, "Mallorca",            38.66025045757143),
    @(18, 17, "Elche",               37.13527028658773),
    @(19, 18, "Levante",             34.28422230441273),
    @(20, 19, "Girona",              32.82276187198359),
    @(21, 20, "Real Oviedo",         29.67982862052838)
)

foreach ($row in $data) {
    $r      = $row[0]
    $rank   = $row[1]
    $team   = $row[2]
    $expPts = $row[3]

    $ws.Cells.Item($r, 1).Value = $rank
    $ws.Cells.Item($r, 2).Value = $team

    # WIN / TOP4 / TOP5 / TOP6 / RELEGATION placeholders — left blank,
    # ready for the upcoming Monte Carlo simulation values.
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = ""
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 7).Value = ""

    # Recalculated ExpPoints now lives in column H.
    $ws.Cells.Item($r, 8).Value = $expPts
}
